# Apply "team member report and sprint backlog" update to Sprint 4 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 4")

# --- Fill in the "Team Member Names and Percentage Contribution" column
# --- for the two existing rows that were missing it.
$ws.Range("G19").Value = "Sakshyam: 100%"
$ws.Range("G20").Value = "Sakshyam: 100%"

# --- Append new backlog rows (22-32) ---
$ws.Range("A22").Value = "144, Update SRS and UC document document for User story 30"
$ws.Range("B22").Value = 30
$ws.Range("C22").Value = "F"
$ws.Range("D22").Value = "T"
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = "Micheal: 100%"

$ws.Range("A23").Value = "148, Update SRS and UC document document for User story 24"
$ws.Range("B23").Value = 24
$ws.Range("C23").Value = "F"
$ws.Range("D23").Value = "T"
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = "Vasilis: 100%"

$ws.Range("A24").Value = "162, Restructure URN document"
$ws.Range("B24").Value = "N/A"
$ws.Range("C24").Value = "N/A"
$ws.Range("D24").Value = "T"
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 3
$ws.Range("G24").Value = "Vasilis: 100%"

$ws.Range("A25").Value = "164,Create use case maps for schedule actions "
$ws.Range("B25").Value = "N/A"
$ws.Range("C25").Value = "N/A"
$ws.Range("D25").Value = "T"
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = "Vasilis: 100%"

$ws.Range("A26").Value = "119, Update Configuration management plan document"
$ws.Range("B26").Value = "N/A"
$ws.Range("C26").Value = "N/A"
$ws.Range("D26").Value = "T"
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = "Vasilis: 100%"

$ws.Range("A27").Value = "114, Update Platform Document"
$ws.Range("B27").Value = "N/A"
$ws.Range("C27").Value = "N/A"
$ws.Range("D27").Value = "T"
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = "Vasilis: 100%"

$ws.Range("A28").Value = "140, update grl and ucm model document for user story #24"
$ws.Range("B28").Value = 24
$ws.Range("C28").Value = "F"
$ws.Range("D28").Value = "T"
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = "Vasilis: 100%"

$ws.Range("A29").Value = "137, update grl and ucm model document for user story #21"
$ws.Range("B29").Value = 21
$ws.Range("C29").Value = "F"
$ws.Range("D29").Value = "T"
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = "Sakshyam: 100%"

$ws.Range("A30").Value = "138, update grl and ucm model document for user story #30"
$ws.Range("B30").Value = 30
$ws.Range("C30").Value = "F"
$ws.Range("D30").Value = "T "
$ws.Range("E30").Value = 2
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = "Micheal: 100%"

$ws.Range("A31").Value = "160, Add unscheduled navigation scenario to URN document"
$ws.Range("B31").Value = "N/A"
$ws.Range("C31").Value = "F"
$ws.Range("D31").Value = "T"
$ws.Range("E31").Value = 2
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = "Vasilis: 100%"

$ws.Range("A32").Value = "161, Add edit course scenario to URN document"
$ws.Range("B32").Value = "N/A"
$ws.Range("C32").Value = "F"
$ws.Range("D32").Value = "T"
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = 2
$ws.Range("G32").Value = "Vasilis: 100%"

# --- Update the view: scroll back to top and select A29 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A29").Select()
